$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = "311-1087-1-ND"
$ws.Range("F7").Value = "311-1087-1-ND"
$ws.Range("F9").Value = "311-1445-1-ND"
$ws.Range("A11").Value = "Q1"
$ws.Range("A12").Value = "Q2"
$ws.Range("F15").Value = "311-1.0KGRCT-ND"
$ws.Range("F16").Value = "311-1.0KGRCT-ND"
$ws.Range("F17").Value = "311-1.0KGRCT-ND"
$ws.Range("F18").Value = "311-1.0KGRCT-ND"
$ws.Range("F24").Value = "311-1.0MGRCT-ND"

$ws.Range("F40").Select()
